# Cost.xlsx: "unify the conception of DataNode, DataTable, Entity."
#
# The sheet that used to model a generic "Property" table is renamed to
# "DataNode" (matching the project-wide rename called out in the commit
# message), and the cursor/selection that Excel persists into the sheet
# view is left sitting on D39 (the last place the author clicked before
# saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")
$ws.Name = "DataNode"

# Re-anchor the saved selection/active cell away from A9 (top of the
# frozen pane) to D39, as captured in the workbook's <selection .../>.
$ws.Range("D39").Select()
